# Weekly update: insert 3 new rows (new price-report date) at the top of
# the data block (row 16) and push the existing rows down by 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 16..18 - this shifts existing rows 16..132 down to 19..135
# and Excel automatically carries the row formatting (incl. the date style on column D).
$ws.Rows("16:18").Insert()

# --- Row 16: new "Especial" entry for 2023-04-19 ---
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 45035
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100107
$ws.Range("H16").Value = "Otros"
$ws.Range("I16").Value = 100107011
$ws.Range("J16").Value = "Tuna"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 16000
$ws.Range("P16").Value = 15500
$ws.Range("Q16").Value = "`$/caja 18 kilos"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 861
$ws.Range("T16").Value = 18

# --- Row 17: new "Primera" entry for 2023-04-19 ---
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 45035
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = "Otros"
$ws.Range("I17").Value = 100107011
$ws.Range("J17").Value = "Tuna"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 13000
$ws.Range("P17").Value = 12500
$ws.Range("Q17").Value = "`$/caja 18 kilos"
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 694
$ws.Range("T17").Value = 18

# --- Row 18: new "Segunda" entry for 2023-04-19 ---
$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 45035
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100107
$ws.Range("H18").Value = "Otros"
$ws.Range("I18").Value = 100107011
$ws.Range("J18").Value = "Tuna"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 9000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 9500
$ws.Range("Q18").Value = "`$/caja 18 kilos"
$ws.Range("R18").Value = "Provincia de Limarí"
$ws.Range("S18").Value = 528
$ws.Range("T18").Value = 18
